# Log_of_all_Blogs.xlsx — add Post61 ("Sleep and Open Terminal using While
# Loop | Shell Scripting") as a new row at the bottom of the Table2 listing
# on Sheet1, then extend the table / dimension to cover it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a fresh row 71 (shifts nothing below it, but makes Excel inherit
# the row-12..70 cell formatting — s="3"/"4"/"5" — instead of leaving the
# new cells unstyled).
$ws.Rows("71:71").Insert(-4121)

# Fill in the new post's data. Touch the link columns before the title so
# the new shared-string entries land in the same order as the workbook's
# own commit (link, title, link).
$ws.Range("E71").Value = "https://programmingport.hashnode.dev/sleep-and-open-terminal-using-while-loop-or-shell-scripting"
$ws.Range("C71").Value = "Sleep and Open Terminal using While Loop | Shell Scripting "
$ws.Range("F71").Value = "https://dev.to/rahulmishra05/sleep-and-open-terminal-using-while-loop-shell-scripting-2mea"
$ws.Range("B71").Value = 61
$ws.Range("D71").Value = "12/14/2020"

# Grow Table2 (and its autofilter) so the new row is part of the table.
$tbl = $ws.ListObjects.Item("Table2")
$tbl.Resize($ws.Range("B10:F71"))

# Match the author's final selection / scroll position.
$win = $excel.ActiveWindow
$win.ScrollColumn = 6
$win.ScrollRow = 51
$ws.Range("F71").Select()
